$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.0155857214680744
$ws.Range("C2").Value = 0.0816155522037875
$ws.Range("D2").Value = 0.601307189542484
$ws.Range("E2").Value = 0.036701860231272
$ws.Range("F2").Value = 0.0583207642031171
$ws.Range("G2").Value = 0.840958605664488
$ws.Range("H2").Value = 0.0717278364337188
$ws.Range("I2").Value = 0.0554717613541143
$ws.Range("J2").Value = 0.772079772079772
$ws.Range("K2").Value = 0.0827886710239651
$ws.Range("L2").Value = 0.0261437908496732
$ws.Range("M2").Value = 0.0110608345902464
$ws.Range("N2").Value = 0.173956762192056
$ws.Range("O2").Value = 0.00251382604323781
$ws.Range("P2").Value = 0.00452488687782805
$ws.Range("Q2").Value = 0.000837942014412603
$ws.Range("R2").Value = 0.955086308027485
$ws.Range("S2").Value = 0.191385956091838
$ws.Range("T2").Value = 0.957935310876487
$ws.Range("U2").Value = 0.0259762024467907
$ws.Range("V2").Value = 0.0888218535277359
$ws.Range("W2").Value = 0.0469247528071058
$ws.Range("X2").Value = 0.0112284229931289

$ws.Range("B3").Value = 0.0884866767219708
$ws.Range("C3").Value = 0.032512150159209
$ws.Range("D3").Value = 0.0315066197419139
$ws.Range("E3").Value = 0.0264789676554382
$ws.Range("F3").Value = 0.841461370873136
$ws.Range("G3").Value = 0.0135746606334842
$ws.Range("H3").Value = 0.0229596111949053
$ws.Range("I3").Value = 0.756326462208815
$ws.Range("J3").Value = 0.0422322775263952
$ws.Range("K3").Value = 0.844310373722138
$ws.Range("L3").Value = 0.00553041729512318
$ws.Range("M3").Value = 0.00955253896430367
$ws.Range("N3").Value = 0.0338528573822691
$ws.Range("O3").Value = 0.000670353611530082
$ws.Range("P3").Value = 0.000335176805765041
$ws.Range("Q3").Value = 0.000335176805765041
$ws.Range("R3").Value = 0.0194402547343724
$ws.Range("S3").Value = 0.0333500921736216
$ws.Range("T3").Value = 0.00620077090665326
$ws.Range("U3").Value = 0.00703871292106586
$ws.Range("V3").Value = 0.011898776604659
$ws.Range("W3").Value = 0.00100553041729512
$ws.Range("X3").Value = 0.00134070722306016

$ws.Range("B4").Value = 0.110440757499581
$ws.Range("C4").Value = 0.0430702195408078
$ws.Range("D4").Value = 0.350092173621585
$ws.Range("E4").Value = 0.156192391486509
$ws.Range("F4").Value = 0.0397184514831574
$ws.Range("G4").Value = 0.0806100217864924
$ws.Range("H4").Value = 0.852187028657617
$ws.Range("I4").Value = 0.0854700854700855
$ws.Range("J4").Value = 0.141277023629965
$ws.Range("K4").Value = 0.0378749790514496
$ws.Range("L4").Value = 0.0221216691804927
$ws.Range("M4").Value = 0.22071392659628
$ws.Range("N4").Value = 0.746941511647394
$ws.Range("O4").Value = 0.00117311882017764
$ws.Range("P4").Value = 0.985587397352103
$ws.Range("Q4").Value = 0.996983408748115
$ws.Range("R4").Value = 0.0212837271660801
$ws.Range("S4").Value = 0.757499581028993
$ws.Range("T4").Value = 0.0326797385620915
$ws.Range("U4").Value = 0.280542986425339
$ws.Range("V4").Value = 0.897938662644545
$ws.Range("W4").Value = 0.934137757667169
$ws.Range("X4").Value = 0.980559745265628

$ws.Range("B5").Value = 0.782637841461371
$ws.Range("C5").Value = 0.841964136081783
$ws.Range("D5").Value = 0.00921736215853863
$ws.Range("E5").Value = 0.77928607340372
$ws.Range("F5").Value = 0.0594938830232948
$ws.Range("G5").Value = 0.0630132394838277
$ws.Range("H5").Value = 0.0516172280878163
$ws.Range("I5").Value = 0.0977040388805095
$ws.Range("J5").Value = 0.0397184514831574
$ws.Range("K5").Value = 0.0336852689793866
$ws.Range("L5").Value = 0.945868945868946
$ws.Range("M5").Value = 0.75532093179152
$ws.Range("N5").Value = 0.0382101558572147
$ws.Range("O5").Value = 0.994972347913524
$ws.Range("P5").Value = 0.00888218535277359
$ws.Range("Q5").Value = 0.00167588402882521
$ws.Range("R5").Value = 0.00318417965476789
$ws.Range("S5").Value = 0.00603318250377074
$ws.Range("T5").Value = 0.00100553041729512
$ws.Range("U5").Value = 0.683928272163566
$ws.Range("V5").Value = 0.00117311882017764
$ws.Range("W5").Value = 0.0177643707055472
$ws.Range("X5").Value = 0.00502765208647562

